$d = $word.ActiveDocument
$d.Content.Find.Execute(", viel besser als die vorherigen Use-Cases", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", viel besser als die vorherigen Use-Cases. Allgemein ist die Arbeit sehr detailliert.", 2)
